{"js": "// PPS: Added VSCode & Jira\n// 1) \"Good Knowledge of ...\" sentence: insert \"Visual Studio Code\" after \"Eclipse, \"\n//    and insert \"Confluence, \" after \"JIRA, \".\n// 2) Skills table, \"Version Control\" row: append \"& Tools\" to the label cell and\n//    \", Jira, Confluence\" to the tools cell.\n\nconst body = context.document.body;\n\n// --- 1) Development tools sentence -------------------------------------\nconst eclipseResults = body.search(\"Eclipse, \", { matchCase: true });\neclipseResults.load(\"text\");\nconst jiraResults = body.search(\"JIRA, Sonar\", { matchCase: true });\njiraResults.load(\"text\");\nawait context.sync();\n\nif (eclipseResults.items.length > 0) {\n  eclipseResults.items[0].insertText(\n    \"Eclipse, Visual Studio Code, \",\n    Word.InsertLocation.replace\n  );\n}\nif (jiraResults.items.length > 0) {\n  jiraResults.items[0].insertText(\n    \"JIRA, Confluence, Sonar\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// --- 2) Skills table \"Version Control\" row ------------------------------\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length > 0) {\n  const table = tables.items[0];\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  let targetRowIndex = -1;\n  // Load all first-column cell texts to find the \"Version Control\" row.\n  const firstCellBodies = [];\n  for (let i = 0; i < rows.items.length; i++) {\n    const c = table.getCell(i, 0);\n    c.body.load(\"text\");\n    firstCellBodies.push(c);\n  }\n  await context.sync();\n\n  for (let i = 0; i < firstCellBodies.length; i++) {\n    if (firstCellBodies[i].body.text.trim() === \"Version Control\") {\n      targetRowIndex = i;\n      break;\n    }\n  }\n\n  if (targetRowIndex >= 0) {\n    const labelCell = table.getCell(targetRowIndex, 0);\n    labelCell.body.insertText(\"& Tools\", Word.InsertLocation.end);\n\n    const toolsCell = table.getCell(targetRowIndex, 1);\n    toolsCell.body.insertText(\", Jira, Confluence\", Word.InsertLocation.end);\n\n    await context.sync();\n  }\n}\n", "ps1": "# PPS: Added VSCode & Jira\n# 1) \"Good Knowledge of ...\" sentence: insert \"Visual Studio Code\" after \"Eclipse, \"\n#    and insert \"Confluence, \" after \"JIRA, \".\n# 2) Skills table, \"Version Control\" row: append \"& Tools\" to the label cell and\n#    \", Jira, Confluence\" to the tools cell.\n\n$d = $word.ActiveDocument\n\n# --- 1) Development tools sentence --------------------------------------\n$range1 = $d.Content\n$found1 = $range1.Find.Execute(\"Eclipse, \")\nif ($found1) {\n    $range1.Text = \"Eclipse, Visual Studio Code, \"\n}\n\n$range2 = $d.Content\n$found2 = $range2.Find.Execute(\"JIRA, Sonar\")\nif ($found2) {\n    $range2.Text = \"JIRA, Confluence, Sonar\"\n}\n\n# --- 2) Skills table \"Version Control\" row -------------------------------\n$table = $d.Tables.Item(1)\n\nfor ($i = 1; $i -le $table.Rows.Count; $i++) {\n    $labelCell = $table.Cell($i, 1)\n    $labelRaw = $labelCell.Range.Text\n    $labelText = $labelRaw.Substring(0, $labelRaw.Length - 2).Trim()\n    if ($labelText -eq \"Version Control\") {\n        $labelRange = $labelCell.Range\n        $labelRange.MoveEnd(1, -2) | Out-Null\n        $labelRange.InsertAfter(\"& Tools\")\n\n        $toolsCell = $table.Cell($i, 2)\n        $toolsRange = $toolsCell.Range\n        $toolsRange.MoveEnd(1, -2) | Out-Null\n        $toolsRange.InsertAfter(\", Jira, Confluence\")\n\n        break\n    }\n}\n\nWrite-Output \"done\"\n"}
